$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.188.32'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.520.64'
$ws.Range('E3').Value = '  +2.02%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '519.23'
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.05'
$ws.Range('E6').Value = '  +0.29%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -0.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.518.16'
$ws.Range('E9').Value = '  +1.83%  '
$ws.Range('E10').Value = '  -1.76%  '
$ws.Range('E11').Value = '  +0.13%  '
$ws.Range('E12').Value = '  -2.73%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.333'
$ws.Range('E13').Value = '  -2.39%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.956.62'
$ws.Range('E14').Value = '  +1.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '58.199.68'
$ws.Range('E15').Value = '  +0.30%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.16'
$ws.Range('E16').Value = '  -0.75%  '
$ws.Range('E17').Value = '  -1.00%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.514.38'
$ws.Range('E18').Value = '  +1.63%  '
$ws.Range('E19').Value = '  -0.25%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.19'
$ws.Range('E20').Value = '  +0.21%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '321.00'
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.05'
$ws.Range('E22').Value = '  +5.23%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.88'
$ws.Range('E24').Value = '  -0.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.404'
$ws.Range('E25').Value = '  -1.76%  '
$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.996'
$ws.Range('E26').Value = '  -0.34%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.161'
$ws.Range('E27').Value = '  +0.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.36'
$ws.Range('E28').Value = '  +0.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0751'
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '168.39'
$ws.Range('E30').Value = '  +1.13%  '
$ws.Range('E31').Value = '  +0.50%  '
$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.25'
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.17'
$ws.Range('E33').Value = '  +0.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.997'
$ws.Range('E35').Value = '  -0.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.06'
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('E37').Value = '  -3.10%  '
$ws.Range('E38').Value = '  -1.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.89'
$ws.Range('E39').Value = '  +0.89%  '
$ws.Range('E40').Value = '  -0.95%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.777'
$ws.Range('E41').Value = '  -1.39%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '277.38'
$ws.Range('E42').Value = '  +2.13%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.44'
$ws.Range('E43').Value = '  -0.83%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.04'
$ws.Range('E44').Value = '  +1.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.598'
$ws.Range('E45').Value = '  +0.87%  '
$ws.Range('E46').Value = '  +1.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '121.90'
$ws.Range('E47').Value = '  -3.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0502'
$ws.Range('E48').Value = '  +2.54%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '17.77'
$ws.Range('E49').Value = '  -0.30%  '
$ws.Range('E50').Value = '  +0.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.99'
$ws.Range('E51').Value = '  -0.09%  '
